# Updated symbol list on Sat Jan 28 19:37:39 UTC 2023 with GitHub Actions
# Refreshes the live crypto price/volume snapshot (cols D/E) and rotates
# three coins (GateToken/KuCoinToken/MXToken/.../LEO, cols B/C) up one row.
# All cells in this sheet are stored as text, so price/percent strings are
# forced back to text (NumberFormat "@" -> Value -> Style "Normal") to avoid
# Excel's automatic "looks like a number" conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue "D2" "306.78"
Set-TextValue "E2" "-0.54%"
Set-TextValue "D3" "38.97"
Set-TextValue "E3" "7.13%"
Set-TextValue "D4" "5.110"
Set-TextValue "E4" "1.06%"
Set-TextValue "D6" "1.952"
Set-TextValue "E6" "-4.24%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D7" "7.967"
Set-TextValue "E7" "1.30%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D8" "0.9310"
Set-TextValue "E8" "0.34%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D9" "0.1466"
Set-TextValue "E9" "2.94%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D10" "0.1929"
Set-TextValue "E10" "-0.15%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D11" "0.09149"
Set-TextValue "E11" "0.15%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D12" "0.03514"
Set-TextValue "E12" "1.86%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D13" "0.09784"
Set-TextValue "E13" "-1.32%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D14" "0.001402"
Set-TextValue "E14" "-0.64%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D15" "0.005873"
Set-TextValue "E15" "-5.69%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D16" "3.789"
Set-TextValue "E16" "-1.25%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D17" "4.178"
Set-TextValue "E17" "0.41%"
Set-TextValue "D18" "3.410"
Set-TextValue "E18" "1.38%"
Set-TextValue "D19" "0.3460"
Set-TextValue "E19" "0.58%"
Set-TextValue "D21" "4.721"
Set-TextValue "E21" "-1.92%"
Set-TextValue "E22" "3.11%"
Set-TextValue "D23" "0.04363"
Set-TextValue "E23" "-0.10%"
Set-TextValue "D24" "0.001237"
Set-TextValue "E24" "0.16%"
Set-TextValue "D25" "0.004284"
Set-TextValue "E25" "-12.91%"
Set-TextValue "E26" "0.04%"
Set-TextValue "D39" "0.02039"
Set-TextValue "E39" "0.48%"
Set-TextValue "D40" "0.05084"
Set-TextValue "E40" "-1.44%"
Set-TextValue "D41" "0.007462"
Set-TextValue "E41" "-0.39%"
Set-TextValue "D42" "0.01029"
Set-TextValue "E42" "1.86%"
Set-TextValue "E43" "-1.65%"
Set-TextValue "E44" "-0.43%"
Set-TextValue "D45" "0.009102"
Set-TextValue "E45" "-6.45%"
Set-TextValue "D46" "0.00006203"
Set-TextValue "E46" "-2.01%"
Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.08%"
Set-TextValue "D48" "0.003101"
Set-TextValue "D50" "0.00002102"
Set-TextValue "E50" "0.08%"
Set-TextValue "D51" "0.0002002"
Set-TextValue "E51" "0.08%"
